$p = $ppt.ActivePresentation
$layout = $p.SlideMaster.CustomLayouts.Item(1)
$sh = $layout.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$sub = $tr.Characters(4,2)
Write-Output "sub: [$($sub.Text)]"
$sub.Text = "30"
$sh2 = $layout.Shapes.Item(3)
Write-Output "AFTER: [$($sh2.TextFrame.TextRange.Text)]"
